$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '39.527.56'
$ws.Cells.Item(2, 5).Value = '  -2.93%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.221.51'
$ws.Cells.Item(3, 5).Value = '  -6.48%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.19%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '297.10'
$ws.Cells.Item(5, 5).Value = '  -4.53%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '82.79'
$ws.Cells.Item(6, 5).Value = '  -4.62%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.511'
$ws.Cells.Item(7, 5).Value = '  -3.75%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.16%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.469'
$ws.Cells.Item(9, 5).Value = '  -4.75%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -7.63%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '29.22'
$ws.Cells.Item(11, 5).Value = '  -3.91%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '47.63'
$ws.Cells.Item(12, 5).Value = '  -9.74%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -2.20%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '2.570.34'
$ws.Cells.Item(14, 5).Value = '  -6.05%  '

# Row 15
$ws.Cells.Item(15, 2).Value = 'Polkadot'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '6.32'
$ws.Cells.Item(15, 5).Value = '  -3.43%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '14.13'
$ws.Cells.Item(16, 5).Value = '  -5.70%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '2.214.03'
$ws.Cells.Item(17, 5).Value = '  -6.89%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.718'
$ws.Cells.Item(18, 5).Value = '  -5.33%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '39.433.96'
$ws.Cells.Item(19, 5).Value = '  -2.90%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '0.0₃0876'
$ws.Cells.Item(20, 5).Value = '  -3.94%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '5.73'
$ws.Cells.Item(21, 5).Value = '  -6.61%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '65.11'
$ws.Cells.Item(22, 5).Value = '  -4.94%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '10.32'
$ws.Cells.Item(23, 5).Value = '  -4.22%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '228.31'
$ws.Cells.Item(24, 5).Value = '  -3.06%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -0.11%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -6.15%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +0.67%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '22.65'
$ws.Cells.Item(28, 5).Value = '  -4.63%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.40%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '9.11'
$ws.Cells.Item(30, 5).Value = '  -1.21%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'Monero'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '149.56'
$ws.Cells.Item(31, 5).Value = '  -2.80%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '32.04'
$ws.Cells.Item(32, 5).Value = '  -6.38%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '1.00'
$ws.Cells.Item(33, 5).Value = '  -0.06%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '4.86'
$ws.Cells.Item(34, 5).Value = '  -6.47%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.0697'
$ws.Cells.Item(35, 5).Value = '  -4.53%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -3.58%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.0971'
$ws.Cells.Item(38, 5).Value = '  -3.14%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '15.26'
$ws.Cells.Item(39, 5).Value = '  -4.10%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '2.64'
$ws.Cells.Item(40, 5).Value = '  -4.96%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '1.64'
$ws.Cells.Item(41, 5).Value = '  -3.39%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '3.65'
$ws.Cells.Item(42, 5).Value = '  -5.06%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '1.910.34'
$ws.Cells.Item(43, 5).Value = '  -2.72%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.0259'
$ws.Cells.Item(44, 5).Value = '  -3.46%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -15.32%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'FraxShare'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '9.05'
$ws.Cells.Item(46, 5).Value = '  -2.86%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'EnergySwap'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '16.20'
$ws.Cells.Item(47, 5).Value = '  -7.94%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  -2.48%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '2.431.01'
$ws.Cells.Item(49, 5).Value = '  -6.38%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '70.71'
$ws.Cells.Item(50, 5).Value = '  -1.66%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '87.35'
$ws.Cells.Item(51, 5).Value = '  -6.37%  '

